# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2210"
#   "<name>_new" -> "<name>_FV2304"
# and wrap the data range in an Excel Table ("Table1") so the header
# row gets a filter + the sheet keeps the header row visible while
# scrolling (frozen pane under row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

# 1. Rename the header row: "_old" suffix becomes "_FV2210",
#    "_new" suffix becomes "_FV2304". Everything else (e.g. "diff")
#    is left untouched.
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -eq $null) { continue }
    $val = [string]$val

    if ($val.EndsWith("_old")) {
        $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2210"
    } elseif ($val.EndsWith("_new")) {
        $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2304"
    }
}

# 2. Freeze the header row so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn the data range into a native Excel Table ("Table1") so the
#    header row exposes the AutoFilter dropdowns.
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$lo = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$lo.Name = "Table1"

Write-Host "Header suffixes updated, table '$($lo.Name)' created over $($dataRange.Address()), header row frozen."
